# Update cryptos list with latest prices / 1h volume change figures
# (GitHub Actions scheduled refresh). Price cells in column D are kept as
# text (leading "'" forces Excel to store them verbatim instead of
# re-parsing as a number), matching how the sheet already stores prices
# as inline strings such as "29.818.98" or "1.009".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.818.98"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "'2.094.01"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'344.61"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.5174"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").Value = "'0.4458"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "'0.09398"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "'52.06"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").Value = "'1.172"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'25.06"
$ws.Range("D13").Value = "'2.098.91"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'6.741"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "'8.036"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "'99.13"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "'0.00001162"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'20.55"
$ws.Range("E19").Value = "  +6.45%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.06703"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'6.166"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").Value = "'29.888.53"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "'12.66"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "'2.318"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "'2.345.63"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "'21.95"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").Value = "'163.84"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "'2.531"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'133.12"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'1.158"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "'1.617"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "'6.218"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").Value = "'3.958"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").Value = "'6.138"
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("D39").Value = "'0.06746"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'0.2275"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'12.46"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").Value = "'0.6880"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'1.296"
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("D44").Value = "'0.6622"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").Value = "'14.18"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("D46").Value = "'2.281"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "'3.640"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "'1.218"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "'0.00000000341"
$ws.Range("E49").Value = "  -7.25%  "
$ws.Range("D50").Value = "'81.66"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'0.07135"
$ws.Range("E51").Value = "  -2.13%  "
